# Weekly update: insert a new "Haba" price record for Macroferia Regional de
# Talca (Maule) as row 14, pushing the existing records (previously rows
# 14-64) down by one row to 15-65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 14 - this shifts rows 14..64 down
# to 15..65 and keeps their contents/formatting intact.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly record.
$ws.Cells.Item(14, 1).Value  = 5
$ws.Cells.Item(14, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(14, 3).Value  = "Maule"
$ws.Cells.Item(14, 4).Value  = 44525
$ws.Cells.Item(14, 5).Value  = 7
$ws.Cells.Item(14, 6).Value  = 100112026
$ws.Cells.Item(14, 7).Value  = "Haba"
$ws.Cells.Item(14, 8).Value  = "Sin especificar"
$ws.Cells.Item(14, 9).Value  = "Primera"
$ws.Cells.Item(14, 10).Value = 300
$ws.Cells.Item(14, 11).Value = 7000
$ws.Cells.Item(14, 12).Value = 7000
$ws.Cells.Item(14, 13).Value = 7000
$ws.Cells.Item(14, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(14, 15).Value = "Región del Maule"
$ws.Cells.Item(14, 16).Value = 280
$ws.Cells.Item(14, 17).Value = 25
$ws.Cells.Item(14, 18).Value = "Hortaliza"
